# allowed camp committee members to generate reports
# - remove the two rejected/duplicate sign-up rows (3 and 7)
# - replace row 11 with an updated committee entry moved down to row 12
# - append newly approved camp-committee report rows (14, 15, 16)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop row 3 (Camp ID 4 / Student 1 / Staff 11) entirely.
$ws.Range("A3:E3").Clear()

# Drop row 7 (Camp ID 7 / Student 8 / Staff 6) entirely.
$ws.Range("A7:E7").Clear()

# The old row 11 entry is superseded - clear it, its replacement lands on row 12.
$ws.Range("A11:E11").Clear()

# Row 12: updated entry (was row 11), now flagged as a camp committee member.
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = 10
$ws.Range("C12").Value = 7
$ws.Range("D12").Value = 2
$ws.Range("E12").Value = $true

# Row 14: newly added participation record.
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = 12
$ws.Range("C14").Value = 14
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = $false

# Row 15: newly added participation record.
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = 2
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = $false

# Row 16: newly added participation record.
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = 11
$ws.Range("C16").Value = 14
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = $false
